$wb = $excel.ActiveWorkbook

# --- "test results" sheet: a new test run overwrote row 11 (questions are now numerated) ---
$ws3 = $wb.Worksheets.Item("test results")
$ws3.Range("B11").Value = "Hitler"
$ws3.Range("C11").Value = "19-01-202500:31:50"
$ws3.Range("E11").Value = 60
$ws3.Range("F11").Value = "'257.32"
$ws3.Range("G11").Value = "(1/6)"
$ws3.Range("H11").Value = "'16.67%"

# the old trailing log rows are gone from the (regenerated) results file
$ws3.Rows("12:14").Delete()

# --- sheet selections / active tab left the way the author last saved them ---
$ws2 = $wb.Worksheets.Item("categories")
[void]$ws2.Activate()
[void]$ws2.Range("B16").Select()

[void]$ws3.Activate()
[void]$ws3.Range("G23").Select()
